$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price cells that are plain decimal numbers stay as text (matching the
# source data format), by forcing Text number format before assigning values.
$r0 = $ws.Range("D5")
$r1 = $ws.Range("D6")
$r2 = $ws.Range("D7")
$r3 = $ws.Range("D8")
$r4 = $ws.Range("D10")
$r5 = $ws.Range("D11")
$r6 = $ws.Range("D15")
$r7 = $ws.Range("D19")
$r8 = $ws.Range("D20")
$r9 = $ws.Range("D21")
$r10 = $ws.Range("D24")
$r11 = $ws.Range("D25")
$r12 = $ws.Range("D26")
$r13 = $ws.Range("D27")
$r14 = $ws.Range("D28")
$r15 = $ws.Range("D31")
$r16 = $ws.Range("D32")
$r17 = $ws.Range("D33")
$r18 = $ws.Range("D34")
$r19 = $ws.Range("D36")
$r20 = $ws.Range("D37")
$r21 = $ws.Range("D38")
$r22 = $ws.Range("D39")
$r23 = $ws.Range("D40")
$r24 = $ws.Range("D42")
$r25 = $ws.Range("D43")
$r26 = $ws.Range("D44")
$r27 = $ws.Range("D45")
$r28 = $ws.Range("D46")
$r29 = $ws.Range("D47")
$r30 = $ws.Range("D48")
$r31 = $ws.Range("D49")
$textRng = $excel.Union($r0, $r1, $r2, $r3, $r4, $r5, $r6, $r7, $r8, $r9, $r10, $r11, $r12, $r13, $r14, $r15, $r16, $r17, $r18, $r19, $r20, $r21, $r22, $r23, $r24, $r25, $r26, $r27, $r28, $r29, $r30, $r31)
foreach ($area in $textRng.Areas) {
    $area.NumberFormat = "@"
}

# Apply updated price (column D) and volume-change (column E) values.
$ws.Range("D2").Value = '63.856.50'
$ws.Range("E2").Value = '  +5.94%  '
$ws.Range("D3").Value = '2.749.52'
$ws.Range("E3").Value = '  +5.28%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '593.45'
$ws.Range("E5").Value = '  +1.61%  '
$ws.Range("D6").Value = '153.12'
$ws.Range("E6").Value = '  +6.98%  '
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("D8").Value = '0.612'
$ws.Range("E8").Value = '  +2.30%  '
$ws.Range("D9").Value = '2.767.80'
$ws.Range("E9").Value = '  +5.62%  '
$ws.Range("D10").Value = '6.73'
$ws.Range("E10").Value = '  +3.48%  '
$ws.Range("D11").Value = '0.114'
$ws.Range("E11").Value = '  +7.88%  '
$ws.Range("E12").Value = '  +3.55%  '
$ws.Range("E13").Value = '  +1.74%  '
$ws.Range("D14").Value = '3.229.18'
$ws.Range("E14").Value = '  +5.10%  '
$ws.Range("D15").Value = '26.63'
$ws.Range("E15").Value = '  +7.24%  '
$ws.Range("D16").Value = '64.076.72'
$ws.Range("E16").Value = '  +6.32%  '
$ws.Range("E17").Value = '  +9.14%  '
$ws.Range("D18").Value = '2.757.59'
$ws.Range("E18").Value = '  +5.47%  '
$ws.Range("D19").Value = '12.09'
$ws.Range("E19").Value = '  +5.73%  '
$ws.Range("D20").Value = '4.92'
$ws.Range("E20").Value = '  +5.99%  '
$ws.Range("D21").Value = '366.05'
$ws.Range("E21").Value = '  +5.60%  '
$ws.Range("E22").Value = '  +1.65%  '
$ws.Range("E23").Value = '  +0.82%  '
$ws.Range("D24").Value = '0.997'
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("D25").Value = '66.28'
$ws.Range("E25").Value = '  +3.98%  '
$ws.Range("D26").Value = '0.169'
$ws.Range("E26").Value = '  +5.42%  '
$ws.Range("D27").Value = '8.68'
$ws.Range("E27").Value = '  +8.39%  '
$ws.Range("D28").Value = '0.996'
$ws.Range("E28").Value = '  -0.30%  '
$ws.Range("D29").Value = '0.0₃0914'
$ws.Range("E29").Value = '  +14.71%  '
$ws.Range("E30").Value = '  +5.00%  '
$ws.Range("D31").Value = '7.13'
$ws.Range("E31").Value = '  +10.04%  '
$ws.Range("D32").Value = '172.20'
$ws.Range("E32").Value = '  +1.78%  '
$ws.Range("D33").Value = '1.21'
$ws.Range("E33").Value = '  +18.39%  '
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").Value = '  -0.09%  '
$ws.Range("E35").Value = '  +6.03%  '
$ws.Range("D36").Value = '4.84'
$ws.Range("E36").Value = '  +12.71%  '
$ws.Range("D37").Value = '1.44'
$ws.Range("E37").Value = '  +9.82%  '
$ws.Range("D38").Value = '1.80'
$ws.Range("E38").Value = '  +9.79%  '
$ws.Range("D39").Value = '1.02'
$ws.Range("E39").Value = '  +19.39%  '
$ws.Range("D40").Value = '349.64'
$ws.Range("E40").Value = '  +9.33%  '
$ws.Range("E41").Value = '  +8.04%  '
$ws.Range("D42").Value = '39.10'
$ws.Range("E42").Value = '  +1.75%  '
$ws.Range("D43").Value = '5.66'
$ws.Range("E43").Value = '  +12.94%  '
$ws.Range("D44").Value = '22.16'
$ws.Range("E44").Value = '  +10.94%  '
$ws.Range("D45").Value = '22.26'
$ws.Range("E45").Value = '  +11.19%  '
$ws.Range("D46").Value = '143.21'
$ws.Range("E46").Value = '  +5.61%  '
$ws.Range("D47").Value = '0.0594'
$ws.Range("E47").Value = '  +7.81%  '
$ws.Range("D48").Value = '0.648'
$ws.Range("E48").Value = '  +6.48%  '
$ws.Range("D49").Value = '0.0259'
$ws.Range("E49").Value = '  +7.54%  '
$ws.Range("E50").Value = '  +2.63%  '
$ws.Range("D51").Value = '2.171.19'
$ws.Range("E51").Value = '  +7.38%  '
